$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("user_id")
$ws2 = $wb.Worksheets.Item("wifi")

$ws1.Range("A5").Value = "Ua1aa40149e47366e6b614488a560c3db"
$ws1.Range("B5").Value = 1
$ws1.Range("A6").Value = "Ufa9b574ba555573176a0f0c217def51f"
$ws1.Range("B6").Value = 1

$ws2.Range("B2").Value = "'45184518"
